$d = $word.ActiveDocument

# The "_GoBack" bookmark currently sits near the end of the document, right
# after "Tim Vigers and Laura Pyle,". It needs to move into the opening
# thank-you paragraph, right after "highlighted in " (splitting off the
# word "yellow" into its own run), and the word "italic " should be removed
# so the sentence reads "...highlighted in yellow." instead of
# "...highlighted in italic yellow."

# Step 1: remove the old (hidden) _GoBack bookmark wherever it currently is.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Step 2: delete the stray word "italic " so the text reads
# "...highlighted in yellow.  " instead of "...highlighted in italic yellow.  "
$find1 = $d.Content
$find1.Find.Execute("highlighted in italic yellow", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "highlighted in yellow", 2) | Out-Null

# Step 3: re-insert the _GoBack bookmark immediately after "highlighted in "
# (i.e. immediately before "yellow."), matching where Word's cursor was left
# after the edit.
$find2 = $d.Content
$find2.Find.Execute("highlighted in ", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0) | Out-Null
$insertionPoint = $d.Range($find2.End, $find2.End)
$d.Bookmarks.Add("_GoBack", $insertionPoint) | Out-Null
